# Weekly timesheet export fix-up: update client names, hours, rate, and totals
# for the week of 2026-01-05 (Doug_Kinsey_2026-01-05.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2026-01-05)
$ws.Range("B2").Value = "Muncey (Maint. Items)"
$ws.Range("C2").Value = 10.5
$ws.Range("E2").Value = 65
$ws.Range("F2").Value = 682.5

# Row 3 (2026-01-06)
$ws.Range("B3").Value = "Ricca"
$ws.Range("C3").Value = 10
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 650

# Row 4 (2026-01-07)
$ws.Range("B4").Value = "Patton"
$ws.Range("C4").Value = 8
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 520

# Row 5 (2026-01-08)
$ws.Range("B5").Value = "Welles"
$ws.Range("C5").Value = 10
$ws.Range("E5").Value = 65
$ws.Range("F5").Value = 650

# Row 6 (2026-01-09, Regular)
$ws.Range("B6").Value = "Caputo"
$ws.Range("C6").Value = 1.5
$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 97.5

# Row 7 (2026-01-09, OT)
$ws.Range("B7").Value = "Caputo"
$ws.Range("E7").Value = 65
$ws.Range("F7").Value = 828.75

# Row 9 (SUBTOTAL)
$ws.Range("F9").Value = 3428.75
